$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.55 = 35967.69 pesos`n✅ 35967.69 pesos = 8.51 = 937.7 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 117
$ws2.Range("O10").Value = 4208.22
$ws2.Range("N12").Value = 4227
$ws2.Range("O12").Value = 110.2
